$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bibi Cell Mundi
$ws.Range("M2").Value = 15255.05
$ws.Range("N2").Value = 15470.5
$ws.Range("AG2").Value = 137509.08

# Row 3 - Bibi Cell Vieiralves
$ws.Range("M3").Value = 6978
$ws.Range("N3").Value = 6459.2
$ws.Range("AG3").Value = 70117.39999999999

# Row 4 - Bibi Cell Manauara
$ws.Range("M4").Value = 2161
$ws.Range("N4").Value = 2235
$ws.Range("AG4").Value = 37906.15

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("M5").Value = 4004.03
$ws.Range("N5").Value = 987
$ws.Range("AG5").Value = 36032.36

# Row 6 - total
$ws.Range("M6").Value = 28398.08
$ws.Range("N6").Value = 25151.7
$ws.Range("AG6").Value = 281564.99
